$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1274.2
$ws.Range("I40").Value = 1183.4445
$ws.Range("J40").Value = 1410.3334
$ws.Range("K40").Value = 1183.4445
$ws.Range("L40").Value = 1410.3334
$ws.Range("M40").Value = -1008.4445
$ws.Range("N40").Value = -1760.3334
$ws.Range("H41").Value = 759
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 898.75
$ws.Range("K41").Value = 200
$ws.Range("L41").Value = 898.75
$ws.Range("M41").Value = 240
$ws.Range("N41").Value = -1778.75
$ws.Range("H94").Value = 2416.1667
$ws.Range("I94").Value = 1899.4
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1899.4
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -1448.4
$ws.Range("N94").Value = -5902
$ws.Range("H96").Value = 31250400
$ws.Range("I96").Value = 35714628
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 107143884
$ws.Range("L96").Value = 2400
$ws.Range("M96").Value = -107142511
$ws.Range("N96").Value = -5146
$ws.Range("H99").Value = 1800
$ws.Range("I99").Value = 900
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 2700
$ws.Range("L99").Value = 6300
$ws.Range("M99").Value = -1202
$ws.Range("N99").Value = -9296
$ws.Range("H100").Value = 2474.1667
$ws.Range("J100").Value = 3741.4285
$ws.Range("L100").Value = 3741.4285
$ws.Range("N100").Value = -4823.4285
$ws.Range("H135").Value = 23817472
$ws.Range("I135").Value = 744.375
$ws.Range("K135").Value = 6699.375
$ws.Range("M135").Value = -4164.375
$ws.Range("H138").Value = 2140.6963
$ws.Range("I138").Value = 1975.8334
$ws.Range("J138").Value = 2189.3442
$ws.Range("K138").Value = 5927.5002
$ws.Range("L138").Value = 6568.0326
$ws.Range("M138").Value = -787.5002000000004
$ws.Range("N138").Value = -16848.0326
$ws.Range("H141").Value = 3933.3333
$ws.Range("I141").Value = 4165
$ws.Range("K141").Value = 12495
$ws.Range("M141").Value = -7315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5575
$ws.Range("I32").Value = 5101.476
$ws.Range("K32").Value = 5101.476
$ws.Range("M32").Value = -4814.476
$ws.Range("H102").Value = 1962.5
$ws.Range("J102").Value = 1962.5
$ws.Range("L102").Value = 1962.5
$ws.Range("N102").Value = -5206.5
$ws.Range("H132").Value = 15655.473
$ws.Range("I132").Value = 1727.5172
$ws.Range("J132").Value = 73357
$ws.Range("K132").Value = 5182.5516
$ws.Range("L132").Value = 220071
$ws.Range("M132").Value = -2652.5516
$ws.Range("N132").Value = -225131

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2546.818
$ws.Range("I94").Value = 1280
$ws.Range("J94").Value = 3602.5
$ws.Range("K94").Value = 1280
$ws.Range("L94").Value = 3602.5
$ws.Range("M94").Value = -829
$ws.Range("N94").Value = -4504.5
$ws.Range("H99").Value = 1505.0526
$ws.Range("I99").Value = 1098.9286
$ws.Range("J99").Value = 2642.2
$ws.Range("K99").Value = 1098.9286
$ws.Range("L99").Value = 2642.2
$ws.Range("M99").Value = 399.0714
$ws.Range("N99").Value = -5638.2
$ws.Range("H105").Value = 3575264.2
$ws.Range("I105").Value = 4212.375
$ws.Range("K105").Value = 4212.375
$ws.Range("M105").Value = -2465.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11789.352
$ws.Range("I31").Value = 16692
$ws.Range("J31").Value = 5354.625
$ws.Range("K31").Value = 16692
$ws.Range("L31").Value = 5354.625
$ws.Range("M31").Value = -16397
$ws.Range("N31").Value = -5944.625
$ws.Range("H34").Value = 11789.352
$ws.Range("I34").Value = 16692
$ws.Range("J34").Value = 5354.625
$ws.Range("K34").Value = 16692
$ws.Range("L34").Value = 5354.625
$ws.Range("M34").Value = -16490
$ws.Range("N34").Value = -5758.625
$ws.Range("H135").Value = 50450
$ws.Range("J135").Value = 50450
$ws.Range("L135").Value = 50450
$ws.Range("N135").Value = -60590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 556324.4399999999
$ws.Range("I129").Value = 1173.3334
$ws.Range("J129").Value = 833900
$ws.Range("K129").Value = 3520.0002
$ws.Range("L129").Value = 2501700
$ws.Range("M129").Value = 1479.9998
$ws.Range("N129").Value = -2511700
$ws.Range("H131").Value = 728.91
$ws.Range("J131").Value = 729
$ws.Range("L131").Value = 2187
$ws.Range("N131").Value = -12267
$ws.Range("H139").Value = 5127.5386
$ws.Range("I139").Value = 2838
$ws.Range("J139").Value = 6558.5
$ws.Range("K139").Value = 8514
$ws.Range("L139").Value = 19675.5
$ws.Range("M139").Value = -3374
$ws.Range("N139").Value = -29955.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1110.8334
$ws.Range("I97").Value = 1043.75
$ws.Range("J97").Value = 1245
$ws.Range("K97").Value = 1043.75
$ws.Range("L97").Value = 1245
$ws.Range("M97").Value = -547.75
$ws.Range("N97").Value = -2237
$ws.Range("H123").Value = 7647.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 876.7143
$ws.Range("I46").Value = 892.4
$ws.Range("J46").Value = 837.5
$ws.Range("K46").Value = 892.4
$ws.Range("L46").Value = 837.5
$ws.Range("M46").Value = -704.4
$ws.Range("N46").Value = -1213.5
$ws.Range("H93").Value = 3336.75
$ws.Range("I93").Value = 2955.7144
$ws.Range("K93").Value = 2955.7144
$ws.Range("M93").Value = -1707.7144
$ws.Range("H100").Value = 2931.818
$ws.Range("I100").Value = 2449.8333
$ws.Range("J100").Value = 3112.5625
$ws.Range("K100").Value = 2449.8333
$ws.Range("L100").Value = 3112.5625
$ws.Range("M100").Value = -1908.8333
$ws.Range("N100").Value = -4194.5625
$ws.Range("H122").Value = 1403653.5
$ws.Range("I122").Value = 2803920.5
$ws.Range("J122").Value = 3386.4285
$ws.Range("K122").Value = 8411761.5
$ws.Range("L122").Value = 10159.2855
$ws.Range("M122").Value = -8409311.5
$ws.Range("N122").Value = -15059.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H96").Value = 2142.8572
$ws.Range("I96").Value = 2350
$ws.Range("J96").Value = 1866.6666
$ws.Range("K96").Value = 2350
$ws.Range("L96").Value = 1866.6666
$ws.Range("M96").Value = -977
$ws.Range("N96").Value = -4612.6666
$ws.Range("H100").Value = 190.72728
$ws.Range("I100").Value = 183.11111
$ws.Range("J100").Value = 225
$ws.Range("K100").Value = 366.22222
$ws.Range("L100").Value = 450
$ws.Range("M100").Value = 174.77778
$ws.Range("N100").Value = -1532
$ws.Range("H126").Value = 1181.1
$ws.Range("J126").Value = 1212.5
$ws.Range("L126").Value = 3637.5
$ws.Range("N126").Value = -8577.5
